# Atualização automática de preços de eletricidade
# Updates the single data row (row 2) of the Spot_PT sheet with the latest
# hourly electricity spot prices, matching the new day's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number, keeps existing YYYY-MM-DD date format)
$ws.Range("A2").Value = 45983

# Hourly prices 0h-1h .. 23h-24h
$ws.Range("B2").Value = 89.52
$ws.Range("C2").Value = 79.16
$ws.Range("D2").Value = 75.45
$ws.Range("E2").Value = 70.03
$ws.Range("F2").Value = 67.22
$ws.Range("G2").Value = 67.22
$ws.Range("H2").Value = 74.47
$ws.Range("I2").Value = 83.77
$ws.Range("J2").Value = 76.94
$ws.Range("K2").Value = 34.68
$ws.Range("L2").Value = 2.13
$ws.Range("M2").Value = 0.53
$ws.Range("N2").Value = 0.65
$ws.Range("O2").Value = 0.6899999999999999
$ws.Range("P2").Value = 0.79
$ws.Range("Q2").Value = 2.78
$ws.Range("R2").Value = 42.91
$ws.Range("S2").Value = 90.25
$ws.Range("T2").Value = 98.84999999999999
$ws.Range("U2").Value = 102.83
$ws.Range("V2").Value = 91.13
$ws.Range("W2").Value = 92.17
$ws.Range("X2").Value = 90.23
$ws.Range("Y2").Value = 83.33

# Daily average price
$ws.Range("Z2").Value = 59.07

# Slot_4h_max (unchanged label) and Slot_4h_price
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 89.22

# Slot_2h_frist (unchanged label) and Slot_2h_frist_price
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 100.84

# Slot_2h_second (unchanged label) and Slot_2h_second_price
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 91.65000000000001

# Slot_min_price window label changes from 2h-16h to 9h-16h
$ws.Range("AG2").Value = "9h-16h"
